# Scrub the Oracle Cloud test credentials (URL / implementation user /
# password) that were hard-coded on the Input_Value sheet, and drop the
# now-stale hyperlink that pointed at the old URL cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# The hyperlink lives on M2 (the URL cell) - remove it before clearing the
# cell it was attached to.
$ws.Hyperlinks.Delete()

# Clear the URL / UserName / Password values - this workbook no longer
# ships with baked-in login credentials.
$ws.Range("M2:O2").ClearContents()
$ws.Rows.Item(1).AutoFit()

# Leave the view focused on the (now blank) credential columns.
$ws.Range("J1").Select()
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("M2:O2").Select()
